# Updated symbol list on Sun Jan 22 23:51:06 UTC 2023 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) figures to the crypto table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. A leading apostrophe is used so Excel
# stores these numeric-looking / percent-looking values as literal text,
# matching the original inline-string cell type in the sheet.
$updates = @{
    'D2' = "'302.05"
    'E2' = "'0.20%"
    'E3' = "'3.00%"
    'D4' = "'4.983"
    'E4' = "'-1.22%"
    'D5' = "'0.07750"
    'E5' = "'0.73%"
    'D6' = "'2.101"
    'E6' = "'-3.48%"
    'D7' = "'7.910"
    'E7' = "'-1.47%"
    'D8' = "'4.037"
    'E8' = "'0.58%"
    'D9' = "'0.9215"
    'E9' = "'-0.87%"
    'D10' = "'0.09785"
    'E10' = "'4.22%"
    'D11' = "'0.1870"
    'E11' = "'2.51%"
    'D12' = "'0.08628"
    'E12' = "'1.71%"
    'D13' = "'0.03512"
    'E13' = "'-2.39%"
    'D14' = "'0.09956"
    'E14' = "'-0.26%"
    'D15' = "'0.001465"
    'E15' = "'-1.34%"
    'D16' = "'0.005663"
    'E16' = "'-1.64%"
    'D17' = "'3.465"
    'E17' = "'-0.39%"
    'D18' = "'2.312"
    'E18' = "'5.84%"
    'D19' = "'0.3409"
    'E19' = "'-1.59%"
    'D20' = "'0.1341"
    'E20' = "'1.02%"
    'D21' = "'4.780"
    'E21' = "'4.30%"
    'D23' = "'0.04591"
    'E23' = "'-1.60%"
    'D24' = "'0.005082"
    'E24' = "'13.61%"
    'D25' = "'0.001228"
    'E25' = "'-0.91%"
    'D26' = "'0.0001397"
    'E26' = "'6.91%"
    'D39' = "'0.01767"
    'E39' = "'2.60%"
    'D40' = "'0.04657"
    'E40' = "'-0.52%"
    'D41' = "'0.007628"
    'E41' = "'-4.01%"
    'D42' = "'0.1391"
    'E42' = "'-0.67%"
    'D43' = "'0.007706"
    'E43' = "'0.20%"
    'D44' = "'0.002235"
    'E44' = "'0.07%"
    'D45' = "'0.01042"
    'E45' = "'16.47%"
    'D46' = "'0.00006193"
    'E46' = "'-0.93%"
    'E47' = "'-0.95%"
    'D48' = "'0.0005789"
    'E48' = "'-0.20%"
    'E49' = "'565.72%"
    'D50' = "'0.001996"
    'E50' = "'-26.37%"
    'D51' = "'0.00002096"
    'E51' = "'-0.95%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
